$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new reminder-row-id column (reuse the existing header
# formatting from the adjacent STATUS column so N1 picks up the same
# bold/bordered style as the rest of row 1)
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "REMINDER_ROW_ID"

# Per-row reminder UUIDs
$ws.Range("N2").Value = "4d7c1477-3f86-45d5-80e3-f0b872efcc6c"
$ws.Range("N3").Value = "b6ea6898-1d75-408e-9d13-7d68e0a9c46a"
$ws.Range("N4").Value = "353e5b8d-8c30-4ada-a163-02abba498742"
$ws.Range("N5").Value = "dc53b0ae-f468-4454-9216-44dd95cf90d6"
$ws.Range("N6").Value = "182d2127-3809-4a2e-bba0-aae36013aa1f"
